$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-26 from 45244 to 45245
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 45245
}

# Row 26 picks up an explicit row height in the target workbook
$ws.Rows.Item(26).RowHeight = 15

# Add new row 27 with data
$ws.Cells.Item(27, 1).Value = "A 57074-2023"
$ws.Cells.Item(27, 2).Value = 45245
$ws.Cells.Item(27, 3).Value = 45245
$ws.Cells.Item(27, 4).Value = "OKÄNT"
$ws.Cells.Item(27, 5).Value = "OKÄNT"
$ws.Cells.Item(27, 6).Value = "SCA"
$ws.Cells.Item(27, 7).Value = 3.2
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 0
$ws.Cells.Item(27, 14).Value = 0
$ws.Cells.Item(27, 15).Value = 0
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(27, 17).Value = 0

# Apply same number format as B and C columns (date style) to B27/C27
$ws.Cells.Item(27, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(27, 3).NumberFormat = "YYYY-MM-DD"

# R27 gets the wrap-text style (same as R2:R26), keep empty value
$ws.Cells.Item(27, 18).WrapText = $true
